{"js": "const replacements = [\n  [\"2025-06-15 Sunday\", \"2025-06-16 Monday\"],\n  [\"15\u00f73=5, 0\", \"40\u00f78=5, 0\"],\n  [\"67\u00f79=7, 4\", \"83\u00f76=13, 5\"],\n  [\"40\u00f72=20, 0\", \"19\u00f72=9, 1\"],\n  [\"11\u00f73=3, 2\", \"77\u00f78=9, 5\"],\n  [\"89\u00f74=22, 1\", \"10\u00f79=1, 1\"],\n  [\"15\u00f77=2, 1\", \"91\u00f72=45, 1\"],\n  [\"87\u00f72=43, 1\", \"26\u00f72=13, 0\"],\n  [\"70\u00f76=11, 4\", \"63\u00f76=10, 3\"],\n  [\"68\u00f74=17, 0\", \"76\u00f78=9, 4\"],\n  [\"34\u00f74=8, 2\", \"24\u00f78=3, 0\"],\n  [\"18\u00f78=2, 2\", \"64\u00f76=10, 4\"],\n  [\"57\u00f77=8, 1\", \"79\u00f72=39, 1\"],\n  [\"99\u00f72=49, 1\", \"89\u00f77=12, 5\"],\n  [\"90\u00f74=22, 2\", \"66\u00f72=33, 0\"],\n  [\"16\u00f79=1, 7\", \"78\u00f75=15, 3\"],\n  [\"21\u00f75=4, 1\", \"14\u00f74=3, 2\"],\n  [\"81\u00f79=9, 0\", \"14\u00f77=2, 0\"],\n  [\"85\u00f76=14, 1\", \"16\u00f76=2, 4\"],\n  [\"47\u00f78=5, 7\", \"59\u00f78=7, 3\"],\n  [\"16\u00f72=8, 0\", \"73\u00f79=8, 1\"],\n  [\"95\u00f72=47, 1\", \"91\u00f75=18, 1\"],\n  [\"62\u00f77=8, 6\", \"39\u00f76=6, 3\"],\n  [\"97\u00f72=48, 1\", \"47\u00f77=6, 5\"],\n  [\"45\u00f72=22, 1\", \"41\u00f73=13, 2\"],\n  [\"93\u00f75=18, 3\", \"76\u00f73=25, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2025-06-15 Sunday\", \"2025-06-16 Monday\")\n    ,@(\"15\u00f73=5, 0\", \"40\u00f78=5, 0\")\n    ,@(\"67\u00f79=7, 4\", \"83\u00f76=13, 5\")\n    ,@(\"40\u00f72=20, 0\", \"19\u00f72=9, 1\")\n    ,@(\"11\u00f73=3, 2\", \"77\u00f78=9, 5\")\n    ,@(\"89\u00f74=22, 1\", \"10\u00f79=1, 1\")\n    ,@(\"15\u00f77=2, 1\", \"91\u00f72=45, 1\")\n    ,@(\"87\u00f72=43, 1\", \"26\u00f72=13, 0\")\n    ,@(\"70\u00f76=11, 4\", \"63\u00f76=10, 3\")\n    ,@(\"68\u00f74=17, 0\", \"76\u00f78=9, 4\")\n    ,@(\"34\u00f74=8, 2\", \"24\u00f78=3, 0\")\n    ,@(\"18\u00f78=2, 2\", \"64\u00f76=10, 4\")\n    ,@(\"57\u00f77=8, 1\", \"79\u00f72=39, 1\")\n    ,@(\"99\u00f72=49, 1\", \"89\u00f77=12, 5\")\n    ,@(\"90\u00f74=22, 2\", \"66\u00f72=33, 0\")\n    ,@(\"16\u00f79=1, 7\", \"78\u00f75=15, 3\")\n    ,@(\"21\u00f75=4, 1\", \"14\u00f74=3, 2\")\n    ,@(\"81\u00f79=9, 0\", \"14\u00f77=2, 0\")\n    ,@(\"85\u00f76=14, 1\", \"16\u00f76=2, 4\")\n    ,@(\"47\u00f78=5, 7\", \"59\u00f78=7, 3\")\n    ,@(\"16\u00f72=8, 0\", \"73\u00f79=8, 1\")\n    ,@(\"95\u00f72=47, 1\", \"91\u00f75=18, 1\")\n    ,@(\"62\u00f77=8, 6\", \"39\u00f76=6, 3\")\n    ,@(\"97\u00f72=48, 1\", \"47\u00f77=6, 5\")\n    ,@(\"45\u00f72=22, 1\", \"41\u00f73=13, 2\")\n    ,@(\"93\u00f75=18, 3\", \"76\u00f73=25, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
